# V86: fixed formatting issues + updated errata list
# Fill in the final (previously blank) row of the errata table with the
# new entry: Violin 2, bars 489/491, "Slur over slur - needs clarification".

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"'

# --- Cell 1: Instrument ---------------------------------------------------
$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$row = $t.Rows.Item($t.Rows.Count)
$cell1 = $row.Cells.Item(1)
$xml1 = '<w:p ' + $wNs + ' w14:paraId="5D2181BB" w14:textId="6CE99B5B" w:rsidR="00990619" w:rsidRDefault="00990619" w:rsidP="00BD52A8"><w:pPr><w:rPr><w:rFonts w:ascii="Times" w:eastAsia="Times" w:hAnsi="Times" w:cs="Times"/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times" w:eastAsia="Times" w:hAnsi="Times" w:cs="Times"/><w:bCs/></w:rPr><w:lastRenderedPageBreak/><w:t>Violin 2</w:t></w:r></w:p>'
$cell1.Range.InsertXML($xml1)

# --- Cell 2: Bar -----------------------------------------------------------
$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$row = $t.Rows.Item($t.Rows.Count)
$cell2 = $row.Cells.Item(2)
$xml2 = '<w:p ' + $wNs + ' w14:paraId="697AC1F1" w14:textId="3AFBE76A" w:rsidR="00990619" w:rsidRDefault="00990619" w:rsidP="00BD52A8"><w:pPr><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Times" w:eastAsia="Times" w:hAnsi="Times" w:cs="Times"/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times" w:eastAsia="Times" w:hAnsi="Times" w:cs="Times"/><w:bCs/></w:rPr><w:t>489, 491</w:t></w:r></w:p>'
$cell2.Range.InsertXML($xml2)

# --- Cell 3: Description ----------------------------------------------------
$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$row = $t.Rows.Item($t.Rows.Count)
$cell3 = $row.Cells.Item(3)
$xml3 = '<w:p ' + $wNs + ' w14:paraId="256F2C25" w14:textId="0457E514" w:rsidR="00990619" w:rsidRDefault="00990619" w:rsidP="00BD52A8"><w:pPr><w:rPr><w:noProof/></w:rPr></w:pPr><w:r><w:rPr><w:noProof/></w:rPr><w:t>Slur over slur ' + [char]0x2013 + ' needs clarification</w:t></w:r></w:p>'
$cell3.Range.InsertXML($xml3)
